$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Foglio1")

# Copy the style of A5 (existing date row) onto A6 so it picks up the date format/border
$ws.Range("A5").Copy()
$ws.Range("A6").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Fill in new sprint day entry
$ws.Range("A6").Value = 43784
$ws.Range("B6").Value = 5

# Update selection to match the recorded cursor position
$ws.Range("A7").Select()
